# chore: update Sheets via scheduled runner
#
# Refreshes the market-price derived columns (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -> H:N) on the
# per-job leve tables. These are plain pulled values (no formulas in this
# workbook), so each touched cell is written directly with its refreshed
# value. A few rows lose their previously-populated LeveProfitNQ/HQ cell
# (H:L collapsed to the NQ-only break-even case) so those are cleared
# instead of zeroed, matching how the source data pull stopped emitting
# them for that row.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2458
$ws.Range("I113").Value = 2346.6667
$ws.Range("J113").Value = 2625
$ws.Range("K113").Value = 2346.6667
$ws.Range("L113").Value = 2625
$ws.Range("M113").Value = 907.3332999999998
$ws.Range("N113").Value = -9133
$ws.Range("H137").Value = 1105.4849
$ws.Range("I137").Value = 1111.2354
$ws.Range("J137").Value = 1099.375
$ws.Range("K137").Value = 3333.7062
$ws.Range("L137").Value = 3298.125
$ws.Range("M137").Value = -783.7062000000001
$ws.Range("N137").Value = -8398.125
$ws.Range("H138").Value = 2353.5908
$ws.Range("J138").Value = 2387.5715
$ws.Range("L138").Value = 7162.7145
$ws.Range("N138").Value = -17442.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1442.7646
$ws.Range("I2").Value = 1404.9231
$ws.Range("J2").Value = 1565.75
$ws.Range("K2").Value = 1404.9231
$ws.Range("L2").Value = 1565.75
$ws.Range("M2").Value = -1291.9231
$ws.Range("N2").Value = -1791.75
$ws.Range("H61").Value = 8335408
$ws.Range("I61").Value = 19609228
$ws.Range("J61").Value = 2584.4348
$ws.Range("K61").Value = 19609228
$ws.Range("L61").Value = 2584.4348
$ws.Range("M61").Value = -19609016
$ws.Range("N61").Value = -3008.4348
$ws.Range("H116").Value = 1442.7646
$ws.Range("I116").Value = 1404.9231
$ws.Range("J116").Value = 1565.75
$ws.Range("K116").Value = 1404.9231
$ws.Range("L116").Value = 1565.75
$ws.Range("M116").Value = 889.0769
$ws.Range("N116").Value = -6153.75
$ws.Range("H122").Value = 1502.3
$ws.Range("I122").Value = 1337.3334
$ws.Range("K122").Value = 4012.0002
$ws.Range("M122").Value = -1562.0002
$ws.Range("H132").Value = 3100.366
$ws.Range("I132").Value = 3174.7036
$ws.Range("J132").Value = 2957
$ws.Range("K132").Value = 9524.110799999999
$ws.Range("L132").Value = 8871
$ws.Range("M132").Value = -6994.110799999999
$ws.Range("N132").Value = -13931
$ws.Range("H136").Value = 8335408
$ws.Range("I136").Value = 19609228
$ws.Range("J136").Value = 2584.4348
$ws.Range("K136").Value = 58827684
$ws.Range("L136").Value = 7753.3044
$ws.Range("M136").Value = -58825134
$ws.Range("N136").Value = -12853.3044

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1442.7646
$ws.Range("I3").Value = 1404.9231
$ws.Range("J3").Value = 1565.75
$ws.Range("K3").Value = 1404.9231
$ws.Range("L3").Value = 1565.75
$ws.Range("M3").Value = -1290.9231
$ws.Range("N3").Value = -1793.75
$ws.Range("H86").Value = 2190
$ws.Range("I86").Value = 2190
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2190
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1067
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2190
$ws.Range("I89").Value = 2190
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10950
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -5334
$ws.Range("N89").ClearContents()
$ws.Range("H94").Value = 966.5
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3823.3674
$ws.Range("I31").Value = 1263.5714
$ws.Range("J31").Value = 5743.2144
$ws.Range("K31").Value = 1263.5714
$ws.Range("L31").Value = 5743.2144
$ws.Range("M31").Value = -968.5714
$ws.Range("N31").Value = -6333.2144
$ws.Range("H34").Value = 3823.3674
$ws.Range("I34").Value = 1263.5714
$ws.Range("J34").Value = 5743.2144
$ws.Range("K34").Value = 1263.5714
$ws.Range("L34").Value = 5743.2144
$ws.Range("M34").Value = -1061.5714
$ws.Range("N34").Value = -6147.2144
$ws.Range("H132").Value = 7248653.5
$ws.Range("I132").Value = 1610.8667
$ws.Range("K132").Value = 4832.6001
$ws.Range("M132").Value = -2302.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 1206.2222
$ws.Range("I45").Value = 800
$ws.Range("J45").Value = 1322.2858
$ws.Range("K45").Value = 2400
$ws.Range("L45").Value = 3966.8574
$ws.Range("M45").Value = -1868
$ws.Range("N45").Value = -5030.857400000001
$ws.Range("H75").Value = 4957
$ws.Range("I75").Value = 3053.25
$ws.Range("J75").Value = 6480
$ws.Range("K75").Value = 9159.75
$ws.Range("L75").Value = 19440
$ws.Range("M75").Value = -8161.75
$ws.Range("N75").Value = -21436
$ws.Range("H78").Value = 4957
$ws.Range("I78").Value = 3053.25
$ws.Range("J78").Value = 6480
$ws.Range("K78").Value = 27479.25
$ws.Range("L78").Value = 58320
$ws.Range("M78").Value = -22487.25
$ws.Range("N78").Value = -68304
$ws.Range("H129").Value = 1174.9
$ws.Range("I129").Value = 539.9167
$ws.Range("J129").Value = 2127.375
$ws.Range("K129").Value = 1619.7501
$ws.Range("L129").Value = 6382.125
$ws.Range("M129").Value = 3380.2499
$ws.Range("N129").Value = -16382.125
$ws.Range("H137").Value = 6886.593
$ws.Range("J137").Value = 4770
$ws.Range("L137").Value = 14310
$ws.Range("N137").Value = -24510

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1272.4615
$ws.Range("I97").Value = 1295.1666
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 1295.1666
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -799.1666
$ws.Range("N97").Value = -1992
$ws.Range("H113").Value = 5915
$ws.Range("I113").Value = 1355.3077
$ws.Range("J113").Value = 11842.6
$ws.Range("K113").Value = 1355.3077
$ws.Range("L113").Value = 11842.6
$ws.Range("M113").Value = 814.6922999999999
$ws.Range("N113").Value = -16182.6
$ws.Range("H122").Value = 4920.4287
$ws.Range("I122").Value = 4726.5
$ws.Range("J122").Value = 4998
$ws.Range("K122").Value = 14179.5
$ws.Range("L122").Value = 14994
$ws.Range("M122").Value = -11729.5
$ws.Range("N122").Value = -19894
$ws.Range("H132").Value = 2021.7906
$ws.Range("I132").Value = 1728.8611
$ws.Range("J132").Value = 3528.2856
$ws.Range("K132").Value = 5186.5833
$ws.Range("L132").Value = 10584.8568
$ws.Range("M132").Value = -2656.5833
$ws.Range("N132").Value = -15644.8568
$ws.Range("H140").Value = 59530
$ws.Range("J140").Value = 59530
$ws.Range("L140").Value = 59530
$ws.Range("N140").Value = -69890

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 71432240
$ws.Range("I40").Value = 333335500
$ws.Range("J40").Value = 4082.7273
$ws.Range("K40").Value = 333335500
$ws.Range("L40").Value = 4082.7273
$ws.Range("M40").Value = -333335364
$ws.Range("N40").Value = -4354.7273
$ws.Range("H122").Value = 6166.6665
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 6166.6665
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 18499.9995
$ws.Range("N122").Value = -23399.9995
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 3314.7715
$ws.Range("I132").Value = 3346.5
$ws.Range("J132").Value = 3281.1765
$ws.Range("K132").Value = 10039.5
$ws.Range("L132").Value = 9843.529500000001
$ws.Range("M132").Value = -7509.5
$ws.Range("N132").Value = -14903.5295

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 100000
$ws.Range("J64").Value = 100000
$ws.Range("L64").Value = 100000
$ws.Range("N64").Value = -100496
$ws.Range("H67").Value = 100000
$ws.Range("J67").Value = 100000
$ws.Range("L67").Value = 100000
$ws.Range("N67").Value = -101716
$ws.Range("H122").Value = 2666.6667
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("N122").Value = -12900.0001
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 5835179.5
$ws.Range("I132").Value = 1662.3704
$ws.Range("J132").Value = 12683221
$ws.Range("K132").Value = 4987.1112
$ws.Range("L132").Value = 38049663
$ws.Range("M132").Value = -2457.1112
$ws.Range("N132").Value = -38054723
